$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -0.2035943559683977
$ws.Range("E2").Value = 0.07658550861757438
$ws.Range("F2").Value = -0.002045344272760061
$ws.Range("G2").Value = 0.01290957909118949
$ws.Range("H2").Value = 0.02004104860421435
$ws.Range("J2").Value = 0.01936369248775929
$ws.Range("K2").Value = -0.05828071817647335
$ws.Range("L2").Value = -0.1776113038751677
$ws.Range("M2").Value = 0.0252683030334444
$ws.Range("N2").Value = 0.05373952863575021
$ws.Range("O2").Value = -0.04572274689720179
$ws.Range("P2").Value = 0.03751865469340517
$ws.Range("Q2").Value = -0.01295695864447187
$ws.Range("R2").Value = -0.03335443267855292
$ws.Range("D3").Value = 0.1841687564732068
$ws.Range("E3").Value = -0.0453368031047253
$ws.Range("F3").Value = 0.01254328074079483
$ws.Range("G3").Value = 0.03946359327769161
$ws.Range("H3").Value = -0.07328817349843601
$ws.Range("J3").Value = -0.06747201087674803
$ws.Range("K3").Value = 0.1117723134284229
$ws.Range("L3").Value = 0.1231112757778762
$ws.Range("M3").Value = -0.05708090279091604
$ws.Range("N3").Value = -0.04163833013386999
$ws.Range("O3").Value = 0.03307880574120719
$ws.Range("P3").Value = 0.002891733744788976
$ws.Range("Q3").Value = 0.05069895762861504
$ws.Range("R3").Value = 0.06532518096092371
$ws.Range("B4").Value = -0.2035943559683977
$ws.Range("C4").Value = 0.1841687564732068
$ws.Range("E4").Value = -0.5509569394695527
$ws.Range("F4").Value = 0.6293828400720779
$ws.Range("G4").Value = 0.1223853848239987
$ws.Range("H4").Value = -0.1110244182334733
$ws.Range("J4").Value = -0.1562932023579501
$ws.Range("K4").Value = 0.2891929810652991
$ws.Range("L4").Value = -0.2935209195746771
$ws.Range("M4").Value = -0.1021273957489789
$ws.Range("N4").Value = -0.3431864906725092
$ws.Range("O4").Value = -0.04354695927498966
$ws.Range("P4").Value = -0.2593700305742481
$ws.Range("Q4").Value = 0.2325284658972369
$ws.Range("R4").Value = -0.3325476846337796
$ws.Range("B5").Value = 0.07658550861757438
$ws.Range("C5").Value = -0.0453368031047253
$ws.Range("D5").Value = -0.5509569394695527
$ws.Range("F5").Value = -0.8853663935137297
$ws.Range("G5").Value = -0.2472242830355185
$ws.Range("H5").Value = 0.2455181015841376
$ws.Range("J5").Value = 0.2476947700305352
$ws.Range("K5").Value = -0.5071212350070473
$ws.Range("L5").Value = 0.5250884604087512
$ws.Range("M5").Value = 0.0290179487843403
$ws.Range("N5").Value = 0.3857564457044987
$ws.Range("O5").Value = -0.01676470645845839
$ws.Range("P5").Value = 0.233130533751236
$ws.Range("Q5").Value = -0.4610320337244821
$ws.Range("R5").Value = 0.425742931537725
$ws.Range("B6").Value = -0.002045344272760061
$ws.Range("C6").Value = 0.01254328074079483
$ws.Range("D6").Value = 0.6293828400720779
$ws.Range("E6").Value = -0.8853663935137297
$ws.Range("G6").Value = 0.2461524043334889
$ws.Range("H6").Value = -0.2573280506804081
$ws.Range("J6").Value = -0.2878351454495143
$ws.Range("K6").Value = 0.3335690869735569
$ws.Range("L6").Value = -0.5305236742909295
$ws.Range("M6").Value = -0.0722913944078757
$ws.Range("N6").Value = -0.4363268788593173
$ws.Range("O6").Value = 0.02368371637615231
$ws.Range("P6").Value = -0.288938292099288
$ws.Range("Q6").Value = 0.2981308511155079
$ws.Range("R6").Value = -0.3980216784125754
$ws.Range("B7").Value = 0.01290957909118949
$ws.Range("C7").Value = 0.03946359327769161
$ws.Range("D7").Value = 0.1223853848239987
$ws.Range("E7").Value = -0.2472242830355185
$ws.Range("F7").Value = 0.2461524043334889
$ws.Range("H7").Value = 0.3242574197365707
$ws.Range("J7").Value = 0.2995101764734562
$ws.Range("K7").Value = 0.3955375268173798
$ws.Range("L7").Value = -0.0132976252817473
$ws.Range("M7").Value = 0.7992225413587013
$ws.Range("N7").Value = 0.1213210703991177
$ws.Range("O7").Value = 0.007483317921326854
$ws.Range("P7").Value = 0.04365159251840433
$ws.Range("Q7").Value = 0.3678978502472589
$ws.Range("R7").Value = 0.008796204539160919
$ws.Range("B8").Value = 0.02004104860421435
$ws.Range("C8").Value = -0.07328817349843601
$ws.Range("D8").Value = -0.1110244182334733
$ws.Range("E8").Value = 0.2455181015841376
$ws.Range("F8").Value = -0.2573280506804081
$ws.Range("G8").Value = 0.3242574197365707
$ws.Range("J8").Value = 0.9593367903415715
$ws.Range("K8").Value = 0.3312575750638658
$ws.Range("L8").Value = 0.3298054067286034
$ws.Range("M8").Value = 0.609254790913381
$ws.Range("N8").Value = 0.09070312376154253
$ws.Range("O8").Value = -0.0625639460328196
$ws.Range("P8").Value = -0.1674921320028981
$ws.Range("Q8").Value = 0.3506151995279669
$ws.Range("R8").Value = 0.0980603607745419
$ws.Range("B10").Value = 0.01936369248775929
$ws.Range("C10").Value = -0.06747201087674803
$ws.Range("D10").Value = -0.1562932023579501
$ws.Range("E10").Value = 0.2476947700305352
$ws.Range("F10").Value = -0.2878351454495143
$ws.Range("G10").Value = 0.2995101764734562
$ws.Range("H10").Value = 0.9593367903415715
$ws.Range("K10").Value = 0.281471712204577
$ws.Range("L10").Value = 0.378371683812899
$ws.Range("M10").Value = 0.5761133745697244
$ws.Range("N10").Value = 0.187578055402155
$ws.Range("O10").Value = -0.01825362543462824
$ws.Range("P10").Value = -0.06572616774404921
$ws.Range("Q10").Value = 0.3201275896166049
$ws.Range("R10").Value = 0.1481057486721723
$ws.Range("B11").Value = -0.05828071817647335
$ws.Range("C11").Value = 0.1117723134284229
$ws.Range("D11").Value = 0.2891929810652991
$ws.Range("E11").Value = -0.5071212350070473
$ws.Range("F11").Value = 0.3335690869735569
$ws.Range("G11").Value = 0.3955375268173798
$ws.Range("H11").Value = 0.3312575750638658
$ws.Range("J11").Value = 0.281471712204577
$ws.Range("L11").Value = -0.226994598066938
$ws.Range("M11").Value = 0.4133306079946346
$ws.Range("N11").Value = -0.1548652599095349
$ws.Range("O11").Value = -0.2846512465783188
$ws.Range("P11").Value = -0.2176722032077225
$ws.Range("Q11").Value = 0.8968577983671087
$ws.Range("R11").Value = -0.1785815821272974
$ws.Range("B12").Value = -0.1776113038751677
$ws.Range("C12").Value = 0.1231112757778762
$ws.Range("D12").Value = -0.2935209195746771
$ws.Range("E12").Value = 0.5250884604087512
$ws.Range("F12").Value = -0.5305236742909295
$ws.Range("G12").Value = -0.0132976252817473
$ws.Range("H12").Value = 0.3298054067286034
$ws.Range("J12").Value = 0.378371683812899
$ws.Range("K12").Value = -0.226994598066938
$ws.Range("M12").Value = 0.1552746577572451
$ws.Range("N12").Value = 0.6253826188684255
$ws.Range("O12").Value = 0.07452806790010968
$ws.Range("P12").Value = 0.4723374501919353
$ws.Range("Q12").Value = -0.219600254447085
$ws.Range("R12").Value = 0.7368945740904715
$ws.Range("B13").Value = 0.0252683030334444
$ws.Range("C13").Value = -0.05708090279091604
$ws.Range("D13").Value = -0.1021273957489789
$ws.Range("E13").Value = 0.0290179487843403
$ws.Range("F13").Value = -0.0722913944078757
$ws.Range("G13").Value = 0.7992225413587013
$ws.Range("H13").Value = 0.609254790913381
$ws.Range("J13").Value = 0.5761133745697244
$ws.Range("K13").Value = 0.4133306079946346
$ws.Range("L13").Value = 0.1552746577572451
$ws.Range("N13").Value = 0.266704296782332
$ws.Range("O13").Value = -0.06402039008662685
$ws.Range("P13").Value = 0.07576955054845577
$ws.Range("Q13").Value = 0.4751217281760241
$ws.Range("R13").Value = 0.1600804616133711
$ws.Range("B14").Value = 0.05373952863575021
$ws.Range("C14").Value = -0.04163833013386999
$ws.Range("D14").Value = -0.3431864906725092
$ws.Range("E14").Value = 0.3857564457044987
$ws.Range("F14").Value = -0.4363268788593173
$ws.Range("G14").Value = 0.1213210703991177
$ws.Range("H14").Value = 0.09070312376154253
$ws.Range("J14").Value = 0.187578055402155
$ws.Range("K14").Value = -0.1548652599095349
$ws.Range("L14").Value = 0.6253826188684255
$ws.Range("M14").Value = 0.266704296782332
$ws.Range("O14").Value = 0.1105593263013635
$ws.Range("P14").Value = 0.9361045980913972
$ws.Range("Q14").Value = -0.1065908832214402
$ws.Range("R14").Value = 0.7273332791172671
$ws.Range("B15").Value = -0.04572274689720179
$ws.Range("C15").Value = 0.03307880574120719
$ws.Range("D15").Value = -0.04354695927498966
$ws.Range("E15").Value = -0.01676470645845839
$ws.Range("F15").Value = 0.02368371637615231
$ws.Range("G15").Value = 0.007483317921326854
$ws.Range("H15").Value = -0.0625639460328196
$ws.Range("J15").Value = -0.01825362543462824
$ws.Range("K15").Value = -0.2846512465783188
$ws.Range("L15").Value = 0.07452806790010968
$ws.Range("M15").Value = -0.06402039008662685
$ws.Range("N15").Value = 0.1105593263013635
$ws.Range("P15").Value = 0.1154152913399076
$ws.Range("Q15").Value = -0.2035253583151278
$ws.Range("R15").Value = -0.04743242383832395
$ws.Range("B16").Value = 0.03751865469340517
$ws.Range("C16").Value = 0.002891733744788976
$ws.Range("D16").Value = -0.2593700305742481
$ws.Range("E16").Value = 0.233130533751236
$ws.Range("F16").Value = -0.288938292099288
$ws.Range("G16").Value = 0.04365159251840433
$ws.Range("H16").Value = -0.1674921320028981
$ws.Range("J16").Value = -0.06572616774404921
$ws.Range("K16").Value = -0.2176722032077225
$ws.Range("L16").Value = 0.4723374501919353
$ws.Range("M16").Value = 0.07576955054845577
$ws.Range("N16").Value = 0.9361045980913972
$ws.Range("O16").Value = 0.1154152913399076
$ws.Range("Q16").Value = -0.1836060800839042
$ws.Range("R16").Value = 0.6826277671728377
$ws.Range("B17").Value = -0.01295695864447187
$ws.Range("C17").Value = 0.05069895762861504
$ws.Range("D17").Value = 0.2325284658972369
$ws.Range("E17").Value = -0.4610320337244821
$ws.Range("F17").Value = 0.2981308511155079
$ws.Range("G17").Value = 0.3678978502472589
$ws.Range("H17").Value = 0.3506151995279669
$ws.Range("J17").Value = 0.3201275896166049
$ws.Range("K17").Value = 0.8968577983671087
$ws.Range("L17").Value = -0.219600254447085
$ws.Range("M17").Value = 0.4751217281760241
$ws.Range("N17").Value = -0.1065908832214402
$ws.Range("O17").Value = -0.2035253583151278
$ws.Range("P17").Value = -0.1836060800839042
$ws.Range("R17").Value = -0.2073451894415046
$ws.Range("B18").Value = -0.03335443267855292
$ws.Range("C18").Value = 0.06532518096092371
$ws.Range("D18").Value = -0.3325476846337796
$ws.Range("E18").Value = 0.425742931537725
$ws.Range("F18").Value = -0.3980216784125754
$ws.Range("G18").Value = 0.008796204539160919
$ws.Range("H18").Value = 0.0980603607745419
$ws.Range("J18").Value = 0.1481057486721723
$ws.Range("K18").Value = -0.1785815821272974
$ws.Range("L18").Value = 0.7368945740904715
$ws.Range("M18").Value = 0.1600804616133711
$ws.Range("N18").Value = 0.7273332791172671
$ws.Range("O18").Value = -0.04743242383832395
$ws.Range("P18").Value = 0.6826277671728377
$ws.Range("Q18").Value = -0.2073451894415046
